# Update column G ("K") values for rows 2, 3, 5, 6.
# This reflects regenerating save_data using K (strikeouts) instead of Strike#.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 3
